$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2859876786626785
$ws.Range("C2").Value = 0.04072312021581581
$ws.Range("E2").Value = 0.1623781245365947
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002434178837647076
$ws.Range("K2").Value = 0.2502009591290175
$ws.Range("M2").Value = 0.2199625890858385
$ws.Range("N2").Value = 1.498895749586309
$ws.Range("O2").Value = 2.489369783813288
$ws.Range("B3").Value = 0.2554847350227192
$ws.Range("C3").Value = 0.03780696162064601
$ws.Range("E3").Value = 0.1509443028518405
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002436387920974495
$ws.Range("K3").Value = 0.2194797068428471
$ws.Range("M3").Value = 0.1985684469320503
$ws.Range("N3").Value = 1.515875854219739
$ws.Range("O3").Value = 2.506379399645269
$ws.Range("B4").Value = 0.2367899142745671
$ws.Range("C4").Value = 0.03600243382960144
$ws.Range("E4").Value = 0.144031037560886
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002437816313680426
$ws.Range("K4").Value = 0.2005965566013117
$ws.Range("M4").Value = 0.185509760858352
$ws.Range("N4").Value = 1.526841926475684
$ws.Range("O4").Value = 2.518268593841782
$ws.Range("B5").Value = 0.2291805798602695
$ws.Range("C5").Value = 0.03526358635185289
$ws.Range("E5").Value = 0.1412405901075644
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002438416555589816
$ws.Range("K5").Value = 0.192896801381039
$ws.Range("M5").Value = 0.1802076871899416
$ws.Range("N5").Value = 1.531446515049218
$ws.Range("O5").Value = 2.523476722304679
$ws.Range("B6").Value = 0.2279176078382363
$ws.Range("C6").Value = 0.03514069156425847
$ws.Range("E6").Value = 0.1407788497103084
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.0024385173236113
$ws.Range("K6").Value = 0.1916179884107265
$ws.Range("M6").Value = 0.1793284559832742
$ws.Range("N6").Value = 1.532219306508812
$ws.Range("O6").Value = 2.524363456951733
$ws.Range("B7").Value = 0.2366872553599535
$ws.Range("C7").Value = 0.03599248355336471
$ws.Range("E7").Value = 0.14399329647906
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002437824335186248
$ws.Range("K7").Value = 0.2004927336419229
$ws.Range("M7").Value = 0.1854381765201651
$ws.Range("N7").Value = 1.526903475640483
$ws.Range("O7").Value = 2.518337362264376
$ws.Range("B8").Value = 0.2754634023556832
$ws.Range("C8").Value = 0.0397205512488128
$ws.Range("E8").Value = 0.1584133931235314
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002434925616761951
$ws.Range("K8").Value = 0.2396126630197841
$ws.Range("M8").Value = 0.2125698068506594
$ws.Range("N8").Value = 1.504638358081538
$ws.Range("O8").Value = 2.494934604476299
$ws.Range("B9").Value = 0.3517610686232047
$ws.Range("C9").Value = 0.04691919054310745
$ws.Range("E9").Value = 0.1875509354059588
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002429810125902649
$ws.Range("K9").Value = 0.3161547087584324
$ws.Range("M9").Value = 0.2663926138133306
$ws.Range("N9").Value = 1.465263066513328
$ws.Range("O9").Value = 2.460520635074545
$ws.Range("B10").Value = 0.4079623710004228
$ws.Range("C10").Value = 0.0521387694426636
$ws.Range("E10").Value = 0.2094979550066967
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002426395139375846
$ws.Range("K10").Value = 0.3722747197723777
$ws.Range("M10").Value = 0.3063219057431326
$ws.Range("N10").Value = 1.438946486429073
$ws.Range("O10").Value = 2.442252415550371
$ws.Range("B11").Value = 0.4335592202214116
$ws.Range("C11").Value = 0.05449808493490593
$ws.Range("E11").Value = 0.2196028706177842
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.00242491539897996
$ws.Range("K11").Value = 0.3977781902158881
$ws.Range("M11").Value = 0.3245727039548569
$ws.Range("N11").Value = 1.427541515991413
$ws.Range("O11").Value = 2.435469167220845
$ws.Range("B12").Value = 0.4432561700751876
$ws.Range("C12").Value = 0.05538930023799082
$ws.Range("E12").Value = 0.2234469832368831
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002424365610376832
$ws.Range("K12").Value = 0.4074316926877657
$ws.Range("M12").Value = 0.331496365059003
$ws.Range("N12").Value = 1.42330425172754
$ws.Range("O12").Value = 2.433120426625351
$ws.Range("B13").Value = 0.4411675881766257
$ws.Range("C13").Value = 0.05519745976182833
$ws.Range("E13").Value = 0.2226182990726002
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.00242448354845054
$ws.Range("K13").Value = 0.4053528272548306
$ws.Range("M13").Value = 0.3300046756634885
$ws.Range("N13").Value = 1.424213194148752
$ws.Range("O13").Value = 2.433616483152377
$ws.Range("B14").Value = 0.4343569167964176
$ws.Range("C14").Value = 0.05457145017255982
$ws.Range("E14").Value = 0.2199187742772253
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002424869956265008
$ws.Range("K14").Value = 0.3985724742879029
$ws.Range("M14").Value = 0.3251420675868957
$ws.Range("N14").Value = 1.427191278860347
$ws.Range("O14").Value = 2.435271525280825
$ws.Range("B15").Value = 0.4301856895556853
$ws.Range("C15").Value = 0.0541877128562902
$ws.Range("E15").Value = 0.2182675349816918
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002425108015449065
$ws.Range("K15").Value = 0.3944187654961127
$ws.Range("M15").Value = 0.3221652050815393
$ws.Range("N15").Value = 1.429026064381407
$ws.Range("O15").Value = 2.436313937618678
$ws.Range("B16").Value = 0.406290136371382
$ws.Range("C16").Value = 0.05198427581191822
$ws.Range("E16").Value = 0.208840027833908
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002426493323648085
$ws.Range("K16").Value = 0.3706074547569074
$ws.Range("M16").Value = 0.3051309193657588
$ws.Range("N16").Value = 1.439703235576303
$ws.Range("O16").Value = 2.44272647406541
$ws.Range("B17").Value = 0.3916385250643373
$ws.Range("C17").Value = 0.05062864664725453
$ws.Range("E17").Value = 0.2030877078185753
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002427362020775339
$ws.Range("K17").Value = 0.3559930948110264
$ws.Range("M17").Value = 0.2947031616056464
$ws.Range("N17").Value = 1.446398501939584
$ws.Range("O17").Value = 2.447051714601088
$ws.Range("B18").Value = 0.3832142106586502
$ws.Range("C18").Value = 0.04984750569771279
$ws.Range("E18").Value = 0.1997905135247535
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002427868617207199
$ws.Range("K18").Value = 0.3475848964683621
$ws.Range("M18").Value = 0.288713557529654
$ws.Range("N18").Value = 1.450302765832237
$ws.Range("O18").Value = 2.449683191980029
$ws.Range("B19").Value = 0.3803623942955028
$ws.Range("C19").Value = 0.049582782174312
$ws.Range("E19").Value = 0.198676090440685
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002428041336630339
$ws.Range("K19").Value = 0.3447376240358153
$ws.Range("M19").Value = 0.2866869829712968
$ws.Range("N19").Value = 1.451633838586583
$ws.Range("O19").Value = 2.450598837597624
$ws.Range("B20").Value = 0.3931979165032828
$ws.Range("C20").Value = 0.0507731027946221
$ws.Range("E20").Value = 0.2036988721423825
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002427268828086933
$ws.Range("K20").Value = 0.3575490701698811
$ws.Range("M20").Value = 0.2958123690048353
$ws.Range("N20").Value = 1.445680260038438
$ws.Range("O20").Value = 2.446576409100885
$ws.Range("B21").Value = 0.4363572714515271
$ws.Range("C21").Value = 0.05475538452328976
$ws.Range("E21").Value = 0.2207112110010812
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002424756172670195
$ws.Range("K21").Value = 0.4005641429789364
$ws.Range("M21").Value = 0.3265699947243803
$ws.Range("N21").Value = 1.426314329297906
$ws.Range("O21").Value = 2.434779427788385
$ws.Range("B22").Value = 0.464587399690231
$ws.Range("C22").Value = 0.05734515969498943
$ws.Range("E22").Value = 0.2319324322752152
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002423175518649183
$ws.Range("K22").Value = 0.4286528489809314
$ws.Range("M22").Value = 0.3467446248446748
$ws.Range("N22").Value = 1.414132951060219
$ws.Range("O22").Value = 2.428351487141327
$ws.Range("B23").Value = 0.4495184840217803
$ws.Range("C23").Value = 0.05596413775444375
$ws.Range("E23").Value = 0.2259339993767924
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002424013531195471
$ws.Range("K23").Value = 0.4136637172819064
$ws.Range("M23").Value = 0.3359703880182323
$ws.Range("N23").Value = 1.420590858380471
$ws.Range("O23").Value = 2.431664778354019
$ws.Range("B24").Value = 0.3924929184298662
$ws.Range("C24").Value = 0.05070779968730221
$ws.Range("E24").Value = 0.2034225339583955
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002427310938292351
$ws.Range("K24").Value = 0.3568456330182528
$ws.Range("M24").Value = 0.2953108793579986
$ws.Range("N24").Value = 1.446004805676379
$ws.Range("O24").Value = 2.446790843573609
$ws.Range("B25").Value = 0.3310941464397388
$ws.Range("C25").Value = 0.04498386430989854
$ws.Range("E25").Value = 0.1795748077034958
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.00243113345597201
$ws.Range("K25").Value = 0.295467566071693
$ws.Range("M25").Value = 0.3300046756634885
$ws.Range("N25").Value = 1.475456311039665
$ws.Range("O25").Value = 2.468599792670972
